$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Воробьева Полина): C7,D7,E7 go from 0 (green-filled style) to 5 (plain bordered style).
# Copy the "no-fill" format from G7 (already style s="2") onto C7:E7, then write the new values.
$ws.Range("G7").Copy()
$ws.Range("C7:E7").PasteSpecial(-4122)
$ws.Range("C7:E7").Value = 5

# Row 16 (Кобзев Богдан): F16 goes from 0 (green-filled style) to 5 (plain bordered style);
# G16 (already plain bordered, empty) gets a value of 5.
$ws.Range("G16").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("F16:G16").Value = 5

# Move the active selection to G16 (frozen-pane top-left normalizes to C4 automatically on save).
$null = $ws.Range("G16").Select()
